$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scaling")

# NFR_code is required on the example row too, so set the "*" placeholder
# value into F2 (under the NFR_code header) alongside the existing
# GNFR_code example in E2.
$ws.Range("F2").Value = "*"

$ws.Range("F3").Select()
